$d = $word.ActiveDocument

$oldText = "(etc.: CSPO, Kinds, Statements LHS, Concepts, RHS)."
$newText = "(etc.: CSPO, Kinds, Statements). ResourceOccurrences LHS, Concepts (ResourceOccurrence Context Kind), RHS:"

# Locate the bullet-list paragraph to edit/extend by its current text,
# before it is modified (so the lookup doesn't depend on a hard-coded
# paragraph index).
$anchor = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd("`r", "`a") -eq $oldText) {
        $anchor = $para
        break
    }
}

# 1. Rewrite the sentence in place.
$anchor.Range.Text = $newText

# 2. Insert nine new bullet-list items right after it. InsertParagraphAfter
#    duplicates the source paragraph's formatting (numPr/pBdr/shd/ind/rPr),
#    matching the existing list style used by its siblings.
$items = @(
    "(Statement, CSPO, Kind)",
    "(Statement, Kind, CSPO)",
    "(Statement, Statement, Statement)",
    "(Kind, Statement, CSPO)",
    "(Kind, CSPO, Statement)",
    "(Kind, Kind, CSPO)",
    "(CSPO, Statement, Kind)",
    "(CSPO, Kind, Statement)",
    "(CSPO, CSPO, Kind)"
)

foreach ($text in $items) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $anchor.Range.InsertBefore($text)
}
